# Fruta / hortaliza, semanal
# Insert two new weekly observation rows (137, 138) for "Arándano (blue)" at
# Mercado Mayorista Lo Valledor de Santiago, pushing the existing rows
# 137:159 down to 139:161.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above the current row 137 (formatting of row above is
# carried into the new rows, same as Excel's default Insert behaviour).
$ws.Rows("137:138").Insert()

# --- Row 137: new "Especial" quality observation -------------------------
$ws.Range("A137").Value = 6
$ws.Range("B137").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C137").Value = "Metropolitana"
$ws.Range("D137").Value = 44474
$ws.Range("E137").Value = 13
$ws.Range("F137").Value = "Fruta"
$ws.Range("G137").Value = 100101
$ws.Range("H137").Value = "Berries"
$ws.Range("I137").Value = 100101001
$ws.Range("J137").Value = "Arándano (blue)"
$ws.Range("K137").Value = "Sin especificar"
$ws.Range("L137").Value = "Especial"
$ws.Range("M137").Value = 250
$ws.Range("N137").Value = 14000
$ws.Range("O137").Value = 14000
$ws.Range("P137").Value = 14000
$ws.Range("Q137").Value = "`$/bandeja 2 kilos"
$ws.Range("R137").Value = "Provincia del Elquí"
$ws.Range("S137").Value = 7000
$ws.Range("T137").Value = 2

# --- Row 138: new "Segunda" quality observation ---------------------------
$ws.Range("A138").Value = 6
$ws.Range("B138").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C138").Value = "Metropolitana"
$ws.Range("D138").Value = 44474
$ws.Range("E138").Value = 13
$ws.Range("F138").Value = "Fruta"
$ws.Range("G138").Value = 100101
$ws.Range("H138").Value = "Berries"
$ws.Range("I138").Value = 100101001
$ws.Range("J138").Value = "Arándano (blue)"
$ws.Range("K138").Value = "Sin especificar"
$ws.Range("L138").Value = "Segunda"
$ws.Range("M138").Value = 1050
$ws.Range("N138").Value = 11000
$ws.Range("O138").Value = 11000
$ws.Range("P138").Value = 11000
$ws.Range("Q138").Value = "`$/bandeja 2 kilos"
$ws.Range("R138").Value = "Región de O'Higgins"
$ws.Range("S138").Value = 5500
$ws.Range("T138").Value = 2
